$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 updates
$ws.Range("I10").Value = 3.4
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("X10").Value = 10
$ws.Range("AG10").Value = 8.5
$ws.Range("AU10").Value = 8.5

# Row 12 updates
$ws.Range("I12").Value = 3.8
$ws.Range("Q12").Value = 1.65
$ws.Range("R12").Value = 2.2
$ws.Range("W12").Value = 9.5
$ws.Range("AA12").Value = 13
$ws.Range("AD12").Value = 7
$ws.Range("AN12").Value = 4
$ws.Range("AU12").Value = 7.5
